$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 45.875
$ws.Range("C2").Value = 37.56000137329102
$ws.Range("D2").Value = 43.04782792388416
